# Fix the title text on the "Conclusion/Recommendations" slide.
# Originally stored as two runs: "CONCLUSION/" + "rECOMMENDATIONS" (the
# second flagged as an autocorrect-ignored spelling error) plus a trailing
# empty endParaRPr run. Replace the whole paragraph text with a single,
# freshly-inserted run containing "CONCLUSION/recommendations" so the
# leftover run-split / endParaRPr goes away.

$p = $ppt.ActivePresentation

foreach ($s in $p.Slides) {
    foreach ($shape in $s.Shapes) {
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "CONCLUSION/rECOMMENDATIONS") {
                $tr.Delete()
                [void]$tr.InsertAfter("CONCLUSION/recommendations")
            }
        }
    }
}
